$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "aada"
$ws.Range("C3").Value = "aadwika"

$ws.Range("A4").Select()
